$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "clave" (password) column values that were corrected/fixed
$ws.Range("B2").Value = "abc12354"
$ws.Range("B3").Value = "xyz78945"
$ws.Range("B4").Value = "pass4546"

# The longer password values widen column B's best-fit width
$ws.Columns("B").ColumnWidth = 8.5

# Move the active selection to B8 (reflecting the last edited location)
$ws.Range("B8").Select()
